$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# --- Row 1 (headers): drop the "key" column and shift every header from
# column C onward one column to the left (into B..N). Column O becomes empty.
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "isMainQuest"
$ws.Range("D1").Value = "isDisplayQuest"
$ws.Range("E1").Value = "prerequisiteKey"
$ws.Range("F1").Value = "description01"
$ws.Range("G1").Value = "actionType01"
$ws.Range("H1").Value = "objectiveKey01"
$ws.Range("I1").Value = "description02"
$ws.Range("J1").Value = "actionType02"
$ws.Range("K1").Value = "objectiveKey02"
$ws.Range("L1").Value = "requiredAmount01"
$ws.Range("M1").Value = "rewardKey"
$ws.Range("N1").Value = "rewardAmount"
$ws.Range("O1").ClearContents()

# --- Row 2 data: the old "key" value (20000) becomes the new "id" value,
# and every subsequent column shifts one to the left.
$ws.Range("A2").Value = 20000
$ws.Range("B2").Value = "An unknown presence.."
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value = "Call my father?"
$ws.Range("G2").Value = "DialogueResponse"
$ws.Range("H2").ClearContents()
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = 1
$ws.Range("O2").ClearContents()

# --- Row 3 data: same left-shift, old "key" (20001) becomes new "id".
$ws.Range("A3").Value = 20001
$ws.Range("B3").Value = "Economic hardship"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 20000
$ws.Range("F3").Value = "Scavenge for coins."
$ws.Range("G3").Value = "Acquire"
$ws.Range("H3").Value = 10304
$ws.Range("I3").ClearContents()
$ws.Range("L3").Value = 2
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = 1
$ws.Range("O3").ClearContents()

# --- Column widths: column A ("key", originally narrow) is gone, so every
# remaining bestFit column width shifts one position to the left.
$ws.Columns("A").ColumnWidth = 6
$ws.Columns("B").ColumnWidth = 20.109375
$ws.Columns("C").ColumnWidth = 11.109375
$ws.Columns("D").ColumnWidth = 12.77734375
$ws.Columns("E").ColumnWidth = 13.6640625
$ws.Columns("F").ColumnWidth = 16.77734375
$ws.Columns("G").ColumnWidth = 16
$ws.Columns("H").ColumnWidth = 13.44140625
$ws.Columns("I").ColumnWidth = 12
$ws.Columns("J").ColumnWidth = 12
$ws.Columns("K").ColumnWidth = 13.44140625
$ws.Columns("L").ColumnWidth = 16.33203125
$ws.Columns("M").ColumnWidth = 9.6640625
$ws.Columns("N").ColumnWidth = 13.33203125
$ws.Columns("O").ColumnWidth = 13.33203125

# --- Update the current selection to match the saved view.
$ws.Range("G10").Select()
